{"js": "// Absenzenliste-Template: widen the \"Name\" (surname) column by 2mm and\n// narrow the \"Vorname\" (first name) column by 2mm (each originally\n// 1418 dxa / 70.9pt -> 1548 dxa / 77.4pt for Name, 1288 dxa / 64.4pt for\n// Vorname). Word treats a table-cell's column width as shared across the\n// whole column, so writing TableCell.columnWidth on one row updates every\n// row's cell in that column plus the <w:tblGrid> entry.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nconst firstRow = rows.items[0];\nconst cells = firstRow.cells;\ncells.load(\"items\");\nawait context.sync();\n\n// Column index 2 = \"Name\" (surname), column index 3 = \"Vorname\" (first name).\ncells.items[2].columnWidth = 1548 / 20; // 77.4pt\ncells.items[3].columnWidth = 1288 / 20; // 64.4pt\n\nawait context.sync();\n", "ps1": "# Absenzenliste-Template: widen the \"Name\" (surname) column by 2mm and\n# narrow the \"Vorname\" (first name) column by 2mm (each originally\n# 1418 dxa / 70.9pt -> 1548 dxa / 77.4pt for Name, 1288 dxa / 64.4pt for\n# Vorname). In Word's object model, setting a cell's Width resizes the\n# whole column (every row's cell in that column plus <w:tblGrid>).\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# Column 3 = \"Name\" (surname), column 4 = \"Vorname\" (first name).\n$table.Cell(1, 3).Width = 1548 / 20   # 77.4pt\n$table.Cell(1, 4).Width = 1288 / 20   # 64.4pt\n"}
